$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

# Copy formatting (styles / number formats) from the last populated row (86)
# down into the three new rows (87-89) so that the new rows match the look
# of existing example rows (text style s=3, date style s=9, etc). Column H is
# skipped on purpose since the source rows have no cell there at all.
$ws.Range("A86:G86").Copy() | Out-Null
$ws.Range("A87:G87").PasteSpecial(-4122) | Out-Null
$ws.Range("A86:G86").Copy() | Out-Null
$ws.Range("A88:G88").PasteSpecial(-4122) | Out-Null
$ws.Range("A86:G86").Copy() | Out-Null
$ws.Range("A89:G89").PasteSpecial(-4122) | Out-Null

$ws.Range("I86:K86").Copy() | Out-Null
$ws.Range("I87:K87").PasteSpecial(-4122) | Out-Null
$ws.Range("I86:K86").Copy() | Out-Null
$ws.Range("I88:K88").PasteSpecial(-4122) | Out-Null
$ws.Range("I86:K86").Copy() | Out-Null
$ws.Range("I89:K89").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 87: E0086 - Using Risk Controls to indicate measures adopted for a specific event
$ws.Range("A87").Value2 = "E0086"
$ws.Range("B87").Value2 = "Using Risk Controls to indicate measures adopted for a specific event"
$ws.Range("C87").Value2 = "<p>This example shows how risk controls can be used to represent the specific measures in place to detect, monitor, address, and inform stakeholders regarding risks, consequences, and impacts. The first part shows a risk with its associated consequence and resulting impact. The second part shows the process where this risk is present, and the specific controls put in place to: (1) monitor the process for vulnerabilities and avoid the source of risk with oversight,; (2) avoid the consequence event; and (3) monitor for occurence of the impact event, and if it occurs to halt it and to inform the affected person.</p>`n<p>In this example the specific details of the control such as the technology being used, or what characteristics are changed when avoiding an event are not specified to reduce the length of content. To express this information, the specific control can reference two events as Risk and Residual Risk respectively - one which shows the initial values such as likelihood, and another one which shows the modified values. Currently the DPVCG does not dictate a specific method for expressing such information, and welcomes suggestions on best practices for the same.</p>"
$ws.Range("D87").Value2 = "E0086.ttl"
$ws.Range("E87").Value2 = "ttl"
$ws.Range("F87").Value2 = "file"
$ws.Range("G87").Value2 = "risk:RiskControl,risk:controls"
$ws.Range("I87").Value2 = "accepted"
$ws.Range("J87").Value2 = 45669
$ws.Range("K87").Value2 = "Harshvardhan J. Pandit"

# --- Row 88: E0087 - Flexibility of RISK taxonomy in expressing varying roles
$ws.Range("A88").Value2 = "E0087"
$ws.Range("B88").Value2 = "Flexibility of RISK taxonomy in expressing varying roles"
$ws.Range("C88").Value2 = "<p>This example shows how the same concept <code>risk:DataBreach</code> can take on different roles across contexts. The first part shows data breach as a consequence of a malware attack, while the second part shows it as a risk source which could lead to identity fraud and misuse. In these, the same concept of data breach being a consequence or a risk source is context-dependant. The RISK extension does not specify or assert that a data breach is always a risk source or risk or consequence - only that it has the potential to be these things. Therefore, inclusion of RISK extension in the graph does not 'pollute' it or cause misinterpretations such as data breach showing up as a consequence even though it has been asserted to (only) be a risk source.</p>`n<p>The third part shows how these different roles are associated with the concept, and how this can be used in UI/UX to provide users with a list of options for selecting risks. It also shows the distinction between marking or annotating a concept to include it in risk identification process (<code>risk:PotentialRisk</code>) and asserting that the concept is applicable as a risk in a process (<code>dpv:hasRisk</code>). This allows use-cases to directly use the RISK extension's subjective classification of concepts in UI/UX - such as to populate the dropdown of risks, and also provides flexibility to ignore that and create their own thesauri or categorisations by using the <code>risk:PotentialRisk</code> concept to annotate concepts, or to create specific subclasses for contextual classifications.</p>"
$ws.Range("D88").Value2 = "E0087.ttl"
$ws.Range("E88").Value2 = "ttl"
$ws.Range("F88").Value2 = "file"
$ws.Range("G88").Value2 = "risk:PotentialRiskSource,risk:PotentialRisk,risk:PotentialConsequence,risk:PotentialImpact"
$ws.Range("I88").Value2 = "accepted"
$ws.Range("J88").Value2 = 45669
$ws.Range("K88").Value2 = "Harshvardhan J. Pandit"

# --- Row 89: E0088 - Expressing impact on specific rights
$ws.Range("A89").Value2 = "E0088"
$ws.Range("B89").Value2 = "Expressing impact on specific rights"
$ws.Range("C89").Value2 = "To express a specific right has been impacted, the relevant rights impact concept is utilised along with a rights impact category from the RISK extension, which together indicate a right is being impacted in the specified manner. In this example, a customer of a company has complained that their GDPR rights have been violated. Upon investigation, it was found that the GDPR's Transparency Right (Article 13) was impacted by limiting the scope of the right as not all processing was represented, and further the right was obstructed as it was not easy to obtain the information on the website."
$ws.Range("D89").Value2 = "E0088.ttl"
$ws.Range("E89").Value2 = "ttl"
$ws.Range("F89").Value2 = "file"
$ws.Range("G89").Value2 = "risk:RightsImpact,eu-gdpr:RightsImpact,risk:RightsViolated,dpv:hasImpact"
$ws.Range("I89").Value2 = "accepted"
$ws.Range("J89").Value2 = 45669
$ws.Range("K89").Value2 = "Harshvardhan J. Pandit"

Write-Output "Added rows 87-89 to Example sheet"
